$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N, shifting the old
# "Late"/Outstanding columns (N,O,P) one slot to the right (-> O,P,Q).
$ws.Columns("N").Insert()

# The newly inserted column keeps the width of the schedule's other
# numeric columns (10 characters), but as an explicit custom width
# rather than an auto-fit one.
$ws.Columns("N").ColumnWidth = 9.23076923076923

# Make the "Repayment Schedule" sheet the active tab/sheet and select
# the cell the author ended up on after the edit.
$ws.Activate()
$ws.Range("R6").Select()
